$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 54 with data, mirroring the structure of existing rows
$ws.Cells.Item(54, 1).Value = 45986
$ws.Cells.Item(54, 2).Value = 2025
$ws.Cells.Item(54, 3).Value = 2.46481303148316
$ws.Cells.Item(54, 4).Value = 2026
$ws.Cells.Item(54, 5).Value = 2.928470412166684

# Copy the style from the cell above (A53) to the new date cell (A54)
$ws.Cells.Item(53, 1).Copy()
$ws.Cells.Item(54, 1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
